$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove consumer rows 12-16 (the data rows).
# Columns A-F and H-I are fully cleared (contents + formatting); column G
# keeps its existing style but loses its value.
$ws.Range("A12:F16").Clear()
$ws.Range("H12:I16").Clear()
$ws.Range("G12:G16").ClearContents()

# Rows 15-16 had taller (14.9pt) auto-fit height because of their text;
# now that they're empty they shrink to the sheet's standard 13.8pt rows.
$ws.Rows("15:16").RowHeight = 13.8

# Update the view: scroll the window back so column A is visible again
# (it had been scrolled right to show column C first) and move the
# selection to B20.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B20").Select()
